# Auto-generated Excel COM-interop edit script
# Applies updated Diabolos_Profits values (currentAveragePrice / Leve price & profit columns)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR worksheets.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 222.5
$ws.Range("I2").Value = 222.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 222.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -109.5
$ws.Range("H15").Value = 965.6591
$ws.Range("I15").Value = 965.6591
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 2896.9773
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -2727.9773
$ws.Range("H64").Value = 3669.2307
$ws.Range("I64").Value = 2450
$ws.Range("J64").Value = 4714.2856
$ws.Range("K64").Value = 2450
$ws.Range("L64").Value = 4714.2856
$ws.Range("M64").Value = -2202
$ws.Range("N64").Value = -5210.2856
$ws.Range("H67").Value = 3669.2307
$ws.Range("I67").Value = 2450
$ws.Range("J67").Value = 4714.2856
$ws.Range("K67").Value = 2450
$ws.Range("L67").Value = 4714.2856
$ws.Range("M67").Value = -1592
$ws.Range("N67").Value = -6430.2856
$ws.Range("H98").Value = 937.4400000000001
$ws.Range("I98").Value = 858.8570999999999
$ws.Range("J98").Value = 1350
$ws.Range("K98").Value = 858.8570999999999
$ws.Range("L98").Value = 1350
$ws.Range("M98").Value = 639.1429000000001
$ws.Range("H111").Value = 21177.646
$ws.Range("I111").Value = 9974.75
$ws.Range("J111").Value = 48064.6
$ws.Range("K111").Value = 29924.25
$ws.Range("L111").Value = 144193.8
$ws.Range("M111").Value = -26857.25
$ws.Range("N111").Value = -150327.8
$ws.Range("H122").Value = 937.4400000000001
$ws.Range("I122").Value = 858.8570999999999
$ws.Range("J122").Value = 1350
$ws.Range("K122").Value = 2576.5713
$ws.Range("L122").Value = 4050
$ws.Range("M122").Value = -126.5712999999996

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 751.79
$ws.Range("I32").Value = 557.6044000000001
$ws.Range("J32").Value = 2715.2222
$ws.Range("K32").Value = 557.6044000000001
$ws.Range("L32").Value = 2715.2222
$ws.Range("M32").Value = -270.6044000000001
$ws.Range("N32").Value = -3289.2222
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H110").Value = 52632956
$ws.Range("I110").Value = 62501324
$ws.Range("J110").Value = 1669.3334
$ws.Range("K110").Value = 62501324
$ws.Range("L110").Value = 1669.3334
$ws.Range("M110").Value = -62499279
$ws.Range("H132").Value = 426759.88
$ws.Range("I132").Value = 251067.88
$ws.Range("J132").Value = 1430714.1
$ws.Range("K132").Value = 753203.64
$ws.Range("L132").Value = 4292142.300000001
$ws.Range("M132").Value = -750673.64

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H107").Value = 8629152
$ws.Range("I107").Value = 5855.9316
$ws.Range("J107").Value = 35730940
$ws.Range("K107").Value = 5855.9316
$ws.Range("L107").Value = 35730940
$ws.Range("M107").Value = -3935.9316
$ws.Range("N107").Value = -35734780
$ws.Range("H134").Value = 1783.8572
$ws.Range("I134").Value = 1464.5
$ws.Range("J134").Value = 3700
$ws.Range("K134").Value = 4393.5
$ws.Range("L134").Value = 11100
$ws.Range("M134").Value = -1858.5
$ws.Range("N134").Value = -16170

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 7500
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 7500
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 7500
$ws.Range("N23").Value = -7980
$ws.Range("H27").Value = 7500
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 7500
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 7500
$ws.Range("N27").Value = -7884
$ws.Range("H58").Value = 1715.1364
$ws.Range("I58").Value = 1105.6364
$ws.Range("J58").Value = 2324.6365
$ws.Range("K58").Value = 1105.6364
$ws.Range("L58").Value = 2324.6365
$ws.Range("M58").Value = -902.6364000000001
$ws.Range("N58").Value = -2730.6365
$ws.Range("H68").Value = 77249.25
$ws.Range("I68").Value = 45000
$ws.Range("J68").Value = 87999
$ws.Range("K68").Value = 45000
$ws.Range("L68").Value = 87999
$ws.Range("M68").Value = -44251
$ws.Range("N68").Value = -89497
$ws.Range("H71").Value = 77249.25
$ws.Range("I71").Value = 45000
$ws.Range("J71").Value = 87999
$ws.Range("K71").Value = 135000
$ws.Range("L71").Value = 263997
$ws.Range("M71").Value = -131256
$ws.Range("N71").Value = -271485
$ws.Range("H74").Value = 48368.5
$ws.Range("I74").Value = 19333.334
$ws.Range("J74").Value = 60812.145
$ws.Range("K74").Value = 19333.334
$ws.Range("L74").Value = 60812.145
$ws.Range("M74").Value = -18459.334
$ws.Range("N74").Value = -62560.145
$ws.Range("H77").Value = 48368.5
$ws.Range("I77").Value = 19333.334
$ws.Range("J77").Value = 60812.145
$ws.Range("K77").Value = 58000.00199999999
$ws.Range("L77").Value = 182436.435
$ws.Range("M77").Value = -53632.00199999999
$ws.Range("N77").Value = -191172.435
$ws.Range("H99").Value = 2710.4614
$ws.Range("I99").Value = 2286.8572
$ws.Range("J99").Value = 3204.6667
$ws.Range("K99").Value = 2286.8572
$ws.Range("L99").Value = 3204.6667
$ws.Range("M99").Value = -788.8571999999999
$ws.Range("H107").Value = 4555.4443
$ws.Range("I107").Value = 5333
$ws.Range("J107").Value = 4166.6665
$ws.Range("K107").Value = 5333
$ws.Range("L107").Value = 4166.6665
$ws.Range("M107").Value = -3413
$ws.Range("N107").Value = -8006.6665
$ws.Range("H108").Value = 285000.25
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 285000.25
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 285000.25
$ws.Range("N108").Value = -292680.25
$ws.Range("H126").Value = 2710.4614
$ws.Range("I126").Value = 2286.8572
$ws.Range("J126").Value = 3204.6667
$ws.Range("K126").Value = 6860.571599999999
$ws.Range("L126").Value = 9614.000100000001
$ws.Range("M126").Value = -4390.571599999999
$ws.Range("H134").Value = 2295.0344
$ws.Range("I134").Value = 1476.5555
$ws.Range("J134").Value = 3634.3635
$ws.Range("K134").Value = 4429.666499999999
$ws.Range("L134").Value = 10903.0905
$ws.Range("M134").Value = -1894.666499999999
$ws.Range("N134").Value = -15973.0905
$ws.Range("H136").Value = 1715.1364
$ws.Range("I136").Value = 1105.6364
$ws.Range("J136").Value = 2324.6365
$ws.Range("K136").Value = 3316.9092
$ws.Range("L136").Value = 6973.9095
$ws.Range("M136").Value = -766.9092000000001
$ws.Range("N136").Value = -12073.9095

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 85.05
$ws.Range("I38").Value = 101.53333
$ws.Range("J38").Value = 35.6
$ws.Range("K38").Value = 304.59999
$ws.Range("L38").Value = 106.8
$ws.Range("M38").Value = 42.40000999999995
$ws.Range("N38").Value = -800.8
$ws.Range("H40").Value = 3624.913
$ws.Range("I40").Value = 56.3
$ws.Range("J40").Value = 6370
$ws.Range("K40").Value = 225.2
$ws.Range("L40").Value = 25480
$ws.Range("M40").Value = -156.2
$ws.Range("N40").Value = -25618
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 9000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -8064
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 27000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -22320
$ws.Range("H132").Value = 963
$ws.Range("I132").Value = 928.3333
$ws.Range("J132").Value = 997.6667
$ws.Range("K132").Value = 8354.9997
$ws.Range("L132").Value = 8979.0003
$ws.Range("M132").Value = -5824.9997
$ws.Range("N132").Value = -14039.0003
$ws.Range("H137").Value = 450
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 450
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 1350
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -11550

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 50000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 50000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 50000
$ws.Range("N15").Value = -50576
$ws.Range("H33").Value = 2000
$ws.Range("I33").Value = 2000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 2000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -1748
$ws.Range("H70").Value = 9201.044
$ws.Range("I70").Value = 9222.588
$ws.Range("J70").Value = 9140
$ws.Range("K70").Value = 9222.588
$ws.Range("L70").Value = 9140
$ws.Range("M70").Value = -8952.588
$ws.Range("H73").Value = 9201.044
$ws.Range("I73").Value = 9222.588
$ws.Range("J73").Value = 9140
$ws.Range("K73").Value = 9222.588
$ws.Range("L73").Value = 9140
$ws.Range("M73").Value = -8286.588
$ws.Range("H81").Value = 50000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 50000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996
$ws.Range("H84").Value = 50000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 50000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984
$ws.Range("H132").Value = 198849.1
$ws.Range("I132").Value = 224417.89
$ws.Range("J132").Value = 7083.1665
$ws.Range("K132").Value = 673253.67
$ws.Range("L132").Value = 21249.4995
$ws.Range("M132").Value = -670723.67
$ws.Range("N132").Value = -26309.4995

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2432.6
$ws.Range("I16").Value = 2020.8
$ws.Range("J16").Value = 3050.3
$ws.Range("K16").Value = 2020.8
$ws.Range("L16").Value = 3050.3
$ws.Range("M16").Value = -1850.8
$ws.Range("N16").Value = -3390.3
$ws.Range("H46").Value = 2687.3125
$ws.Range("I46").Value = 1499.3334
$ws.Range("J46").Value = 2961.4614
$ws.Range("K46").Value = 1499.3334
$ws.Range("L46").Value = 2961.4614
$ws.Range("M46").Value = -1311.3334
$ws.Range("N46").Value = -3337.4614
$ws.Range("H61").Value = 6817.7
$ws.Range("I61").Value = 7607.1665
$ws.Range("J61").Value = 3659.8333
$ws.Range("K61").Value = 7607.1665
$ws.Range("L61").Value = 3659.8333
$ws.Range("M61").Value = -7405.1665
$ws.Range("H113").Value = 6817.7
$ws.Range("I113").Value = 7607.1665
$ws.Range("J113").Value = 3659.8333
$ws.Range("K113").Value = 7607.1665
$ws.Range("L113").Value = 3659.8333
$ws.Range("M113").Value = -5437.1665

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2651.0754
$ws.Range("I136").Value = 1823.6428
$ws.Range("J136").Value = 5810.364
$ws.Range("K136").Value = 5470.928400000001
$ws.Range("L136").Value = 17431.092
$ws.Range("M136").Value = -2920.928400000001
$ws.Range("N136").Value = -22531.092

Write-Host "Applied Diabolos_Profits updates across all worksheets."
